$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 21 de Marzo de 2020 a las 00:46'

$ws.Cells.Item(9, 2).Value = 19382
$ws.Cells.Item(9, 3).Value = 5593
$ws.Cells.Item(9, 5).Value = 18979
$ws.Cells.Item(9, 7).Value = 49
$ws.Cells.Item(9, 8).Value = 256

$ws.Cells.Item(23, 1).Value = 'Brasil'
$ws.Cells.Item(23, 2).Value = 970
$ws.Cells.Item(23, 3).Value = 330
$ws.Cells.Item(23, 4).Value = 2
$ws.Cells.Item(23, 5).Value = 957
$ws.Cells.Item(23, 6).Value = 18
$ws.Cells.Item(23, 7).Value = 4
$ws.Cells.Item(23, 8).Value = 11

$ws.Cells.Item(24, 1).Value = 'Japon'
$ws.Cells.Item(24, 2).Value = 963
$ws.Cells.Item(24, 3).Value = 20
$ws.Cells.Item(24, 4).Value = 215
$ws.Cells.Item(24, 5).Value = 715
$ws.Cells.Item(24, 6).Value = 50
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(24, 8).Value = 33

$ws.Cells.Item(56, 1).Value = 'Panama'
$ws.Cells.Item(56, 2).Value = 200
$ws.Cells.Item(56, 3).Value = 63
$ws.Cells.Item(56, 4).Value = 1
$ws.Cells.Item(56, 5).Value = 198
$ws.Cells.Item(56, 6).Value = 7
$ws.Cells.Item(56, 8).Value = 1

$ws.Cells.Item(57, 1).Value = 'Libano'
$ws.Cells.Item(57, 2).Value = 177
$ws.Cells.Item(57, 3).Value = 20
$ws.Cells.Item(57, 5).Value = 169
$ws.Cells.Item(57, 6).Value = 3
$ws.Cells.Item(57, 8).Value = 4

$ws.Cells.Item(58, 1).Value = 'Mexico'
$ws.Cells.Item(58, 2).Value = 164
$ws.Cells.Item(58, 3).Value = 46
$ws.Cells.Item(58, 4).Value = 4
$ws.Cells.Item(58, 5).Value = 159
$ws.Cells.Item(58, 6).Value = 1
$ws.Cells.Item(58, 8).Value = 1

$ws.Cells.Item(59, 1).Value = 'Kuwait'
$ws.Cells.Item(59, 2).Value = 159
$ws.Cells.Item(59, 3).Value = 11
$ws.Cells.Item(59, 4).Value = 22
$ws.Cells.Item(59, 5).Value = 137
$ws.Cells.Item(59, 6).Value = 5

$ws.Cells.Item(60, 1).Value = 'Argentina'
$ws.Cells.Item(60, 2).Value = 158
$ws.Cells.Item(60, 3).Value = 30
$ws.Cells.Item(60, 4).Value = 3
$ws.Cells.Item(60, 5).Value = 152
$ws.Cells.Item(60, 6).Value = 0
$ws.Cells.Item(60, 8).Value = 3

$ws.Cells.Item(61, 1).Value = 'Colombia'
$ws.Cells.Item(61, 2).Value = 145
$ws.Cells.Item(61, 3).Value = 37
$ws.Cells.Item(61, 4).Value = 1
$ws.Cells.Item(61, 5).Value = 144
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(61, 8).Value = 0

$ws.Cells.Item(62, 1).Value = 'San Marino'
$ws.Cells.Item(62, 2).Value = 144
$ws.Cells.Item(62, 3).Value = 0
$ws.Cells.Item(62, 4).Value = 4
$ws.Cells.Item(62, 5).Value = 126
$ws.Cells.Item(62, 6).Value = 12
$ws.Cells.Item(62, 8).Value = 14

$ws.Cells.Item(63, 1).Value = 'Emiratos Arabes Unidos'
$ws.Cells.Item(63, 2).Value = 140
$ws.Cells.Item(63, 4).Value = 31
$ws.Cells.Item(63, 5).Value = 107
$ws.Cells.Item(63, 6).Value = 2
$ws.Cells.Item(63, 7).Value = 2
$ws.Cells.Item(63, 8).Value = 2

$ws.Cells.Item(64, 1).Value = 'Eslovaquia'
$ws.Cells.Item(64, 2).Value = 137
$ws.Cells.Item(64, 3).Value = 13
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 137

$ws.Cells.Item(65, 1).Value = 'Armenia'
$ws.Cells.Item(65, 2).Value = 136
$ws.Cells.Item(65, 3).Value = 14
$ws.Cells.Item(65, 4).Value = 1
$ws.Cells.Item(65, 5).Value = 135
$ws.Cells.Item(65, 6).Value = 2
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = 0

$ws.Cells.Item(66, 1).Value = 'Serbia'
$ws.Cells.Item(66, 3).Value = 32
$ws.Cells.Item(66, 4).Value = 2
$ws.Cells.Item(66, 5).Value = 132
$ws.Cells.Item(66, 6).Value = 4
$ws.Cells.Item(66, 8).Value = 1

$ws.Cells.Item(67, 1).Value = 'Taiwan'
$ws.Cells.Item(67, 2).Value = 135
$ws.Cells.Item(67, 3).Value = 27
$ws.Cells.Item(67, 4).Value = 28
$ws.Cells.Item(67, 5).Value = 105
$ws.Cells.Item(67, 7).Value = 1
$ws.Cells.Item(67, 8).Value = 2

$ws.Cells.Item(68, 1).Value = 'Croacia'
$ws.Cells.Item(68, 2).Value = 130
$ws.Cells.Item(68, 3).Value = 20
$ws.Cells.Item(68, 4).Value = 5
$ws.Cells.Item(68, 5).Value = 124
$ws.Cells.Item(68, 8).Value = 1

$ws.Cells.Item(72, 2).Value = 110
$ws.Cells.Item(72, 3).Value = 31
$ws.Cells.Item(72, 5).Value = 110

$ws.Cells.Item(89, 1).Value = 'Venezuela'
$ws.Cells.Item(89, 2).Value = 65
$ws.Cells.Item(89, 3).Value = 23
$ws.Cells.Item(89, 4).Value = 1
$ws.Cells.Item(89, 5).Value = 64
$ws.Cells.Item(89, 6).Value = 0

$ws.Cells.Item(90, 1).Value = 'Malta'
$ws.Cells.Item(90, 2).Value = 64
$ws.Cells.Item(90, 3).Value = 11
$ws.Cells.Item(90, 4).Value = 2

$ws.Cells.Item(91, 1).Value = 'Lituania'
$ws.Cells.Item(91, 2).Value = 63
$ws.Cells.Item(91, 5).Value = 62
$ws.Cells.Item(91, 6).Value = 1
$ws.Cells.Item(91, 8).Value = 0

$ws.Cells.Item(92, 1).Value = 'Tunez'
$ws.Cells.Item(92, 2).Value = 54
$ws.Cells.Item(92, 3).Value = 15
$ws.Cells.Item(92, 4).Value = 1
$ws.Cells.Item(92, 5).Value = 52
$ws.Cells.Item(92, 6).Value = 7
$ws.Cells.Item(92, 8).Value = 1

$ws.Cells.Item(93, 1).Value = 'Nueva Zelanda'
$ws.Cells.Item(93, 2).Value = 53
$ws.Cells.Item(93, 3).Value = 25
$ws.Cells.Item(93, 5).Value = 53

$ws.Cells.Item(94, 1).Value = 'Kazajistan'
$ws.Cells.Item(94, 2).Value = 52
$ws.Cells.Item(94, 3).Value = 8
$ws.Cells.Item(94, 4).Value = 0
$ws.Cells.Item(94, 5).Value = 52

$ws.Cells.Item(95, 1).Value = 'Camboya'
$ws.Cells.Item(95, 2).Value = 51
$ws.Cells.Item(95, 3).Value = 14
$ws.Cells.Item(95, 4).Value = 1
$ws.Cells.Item(95, 5).Value = 50

$ws.Cells.Item(96, 1).Value = 'Oman'
$ws.Cells.Item(96, 3).Value = 0
$ws.Cells.Item(96, 4).Value = 13
$ws.Cells.Item(96, 5).Value = 35

$ws.Cells.Item(97, 1).Value = 'Estado de Palestina'
$ws.Cells.Item(97, 2).Value = 48
$ws.Cells.Item(97, 3).Value = 1
$ws.Cells.Item(97, 4).Value = 17
$ws.Cells.Item(97, 5).Value = 31

$ws.Cells.Item(98, 1).Value = 'Senegal'
$ws.Cells.Item(98, 2).Value = 47
$ws.Cells.Item(98, 3).Value = 11
$ws.Cells.Item(98, 4).Value = 5
$ws.Cells.Item(98, 5).Value = 42

$ws.Cells.Item(99, 1).Value = 'Guadalupe'
$ws.Cells.Item(99, 2).Value = 45
$ws.Cells.Item(99, 3).Value = 12
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 45
$ws.Cells.Item(99, 6).Value = 0

$ws.Cells.Item(100, 1).Value = 'Georgia'
$ws.Cells.Item(100, 3).Value = 4
$ws.Cells.Item(100, 4).Value = 1
$ws.Cells.Item(100, 5).Value = 43
$ws.Cells.Item(100, 6).Value = 1
$ws.Cells.Item(100, 8).Value = 0

$ws.Cells.Item(101, 1).Value = 'Azerbaiyan'
$ws.Cells.Item(101, 2).Value = 44
$ws.Cells.Item(101, 4).Value = 7
$ws.Cells.Item(101, 5).Value = 36
$ws.Cells.Item(101, 8).Value = 1

$ws.Cells.Item(112, 1).Value = 'Jamaica'
$ws.Cells.Item(112, 2).Value = 19
$ws.Cells.Item(112, 4).Value = 2
$ws.Cells.Item(112, 5).Value = 16
$ws.Cells.Item(112, 8).Value = 1

$ws.Cells.Item(113, 1).Value = 'Paraguay'
$ws.Cells.Item(113, 2).Value = 18
$ws.Cells.Item(113, 3).Value = 5
$ws.Cells.Item(113, 5).Value = 18
$ws.Cells.Item(113, 6).Value = 1

$ws.Cells.Item(114, 1).Value = 'Consejo Danes para los Refugiados'
$ws.Cells.Item(114, 2).Value = 18
$ws.Cells.Item(114, 3).Value = 4
$ws.Cells.Item(114, 4).Value = 0
$ws.Cells.Item(114, 5).Value = 18

$ws.Cells.Item(115, 1).Value = 'Ruanda'
$ws.Cells.Item(115, 2).Value = 17
$ws.Cells.Item(115, 3).Value = 6
$ws.Cells.Item(115, 5).Value = 17

$ws.Cells.Item(116, 1).Value = 'Macao'
$ws.Cells.Item(116, 2).Value = 17
$ws.Cells.Item(116, 3).Value = 0
$ws.Cells.Item(116, 4).Value = 10
$ws.Cells.Item(116, 5).Value = 7

$ws.Cells.Item(117, 1).Value = 'Ghana'
$ws.Cells.Item(117, 5).Value = 16
$ws.Cells.Item(117, 8).Value = 0

$ws.Cells.Item(118, 1).Value = 'Bolivia'
$ws.Cells.Item(118, 4).Value = 0
$ws.Cells.Item(118, 5).Value = 16
$ws.Cells.Item(118, 8).Value = 0

$ws.Cells.Item(119, 1).Value = 'Cuba'
$ws.Cells.Item(119, 2).Value = 16
$ws.Cells.Item(119, 3).Value = 5
$ws.Cells.Item(119, 8).Value = 1

$ws.Cells.Item(120, 1).Value = 'Guayana Francesa'
$ws.Cells.Item(120, 3).Value = 0
$ws.Cells.Item(120, 5).Value = 15
$ws.Cells.Item(120, 8).Value = 0

$ws.Cells.Item(121, 1).Value = 'Guyana'
$ws.Cells.Item(121, 2).Value = 15
$ws.Cells.Item(121, 3).Value = 10
$ws.Cells.Item(121, 8).Value = 1

$ws.Cells.Item(122, 1).Value = 'Guam'
$ws.Cells.Item(122, 3).Value = 2

$ws.Cells.Item(123, 1).Value = 'Montenegro'
$ws.Cells.Item(123, 3).Value = 1

$ws.Cells.Item(124, 1).Value = 'Puerto Rico'
$ws.Cells.Item(124, 2).Value = 14
$ws.Cells.Item(124, 3).Value = 8
$ws.Cells.Item(124, 5).Value = 14
$ws.Cells.Item(124, 6).Value = 0

$ws.Cells.Item(129, 1).Value = 'Polinesia Francesa'
$ws.Cells.Item(129, 3).Value = 5

$ws.Cells.Item(130, 1).Value = 'Monaco'
$ws.Cells.Item(130, 3).Value = 1

$ws.Cells.Item(132, 1).Value = 'Trinidad yTobago'
$ws.Cells.Item(132, 3).Value = 0

$ws.Cells.Item(134, 1).Value = 'Etiopia'
$ws.Cells.Item(134, 3).Value = 2

$ws.Cells.Item(138, 1).Value = 'Mayotte'
$ws.Cells.Item(138, 3).Value = 2

$ws.Cells.Item(139, 1).Value = 'Barbados'
$ws.Cells.Item(139, 3).Value = 1

$ws.Cells.Item(140, 1).Value = 'Kirguistan'
$ws.Cells.Item(140, 3).Value = 3

$ws.Cells.Item(141, 1).Value = 'Tanzania'
$ws.Cells.Item(141, 3).Value = 0

$ws.Cells.Item(142, 1).Value = 'Guinea Ecuatorial'

$ws.Cells.Item(143, 1).Value = 'Mongolia'

$ws.Cells.Item(145, 1).Value = 'Bahamas'

$ws.Cells.Item(147, 1).Value = 'San Martin (Parte Francesa)'

$ws.Cells.Item(149, 1).Value = 'Namibia'

$ws.Cells.Item(150, 1).Value = 'San Bartolome'
$ws.Cells.Item(150, 3).Value = 0

$ws.Cells.Item(151, 1).Value = 'Republica de Africa Central'
$ws.Cells.Item(151, 3).Value = 2

$ws.Cells.Item(154, 1).Value = 'Islas Virgenes de los Estados Unidos'

$ws.Cells.Item(157, 1).Value = 'Butan'

$ws.Cells.Item(158, 1).Value = 'Liberia'

$ws.Cells.Item(159, 1).Value = 'Isla de Man'
$ws.Cells.Item(159, 3).Value = 1

$ws.Cells.Item(160, 1).Value = 'Nicaragua'
$ws.Cells.Item(160, 3).Value = 1

$ws.Cells.Item(161, 1).Value = 'Santa Lucia'
$ws.Cells.Item(161, 3).Value = 0

$ws.Cells.Item(162, 1).Value = 'Benin'

$ws.Cells.Item(163, 1).Value = 'Groenlandia'

$ws.Cells.Item(164, 1).Value = 'Nueva Caledonia'

$ws.Cells.Item(165, 1).Value = 'Zambia'
$ws.Cells.Item(165, 3).Value = 0

$ws.Cells.Item(166, 1).Value = 'Mauritania'

$ws.Cells.Item(167, 1).Value = 'Haiti'
$ws.Cells.Item(167, 3).Value = 2

$ws.Cells.Item(168, 1).Value = 'Bermudas'
$ws.Cells.Item(168, 3).Value = 0

$ws.Cells.Item(169, 1).Value = 'Guinea'
$ws.Cells.Item(169, 3).Value = 1

$ws.Cells.Item(171, 1).Value = 'Papua Nueva Guinea'
$ws.Cells.Item(171, 3).Value = 1

$ws.Cells.Item(172, 1).Value = 'Republica de Yibuti'

$ws.Cells.Item(173, 1).Value = 'Fiyi'
$ws.Cells.Item(173, 3).Value = 0

$ws.Cells.Item(174, 1).Value = 'Zimbabue'
$ws.Cells.Item(174, 3).Value = 1

$ws.Cells.Item(175, 1).Value = 'San Vicente y las Granadinas'

$ws.Cells.Item(176, 1).Value = 'Gambia'

$ws.Cells.Item(177, 1).Value = 'Montserrat'
$ws.Cells.Item(177, 3).Value = 0

$ws.Cells.Item(178, 1).Value = 'San Martin (Parte Holandesa)'

$ws.Cells.Item(179, 1).Value = 'El Salvador'

$ws.Cells.Item(180, 1).Value = 'Suazilandia'
$ws.Cells.Item(180, 3).Value = 0

$ws.Cells.Item(181, 1).Value = 'Republica del Chad'

$ws.Cells.Item(182, 1).Value = 'Niger'

$ws.Cells.Item(183, 1).Value = 'Santa Sede'

$ws.Cells.Item(184, 1).Value = 'Somalia'
$ws.Cells.Item(184, 3).Value = 0

$ws.Cells.Item(185, 1).Value = 'Cabo Verde'
$ws.Cells.Item(185, 3).Value = 1

$ws.Cells.Item(186, 1).Value = 'Angola'
$ws.Cells.Item(186, 3).Value = 1

$ws.Cells.Item(187, 1).Value = 'Antigua y Barbuda'
